# Apply the "Final version by Maor" edit to the "The code" section of the document.
#
# 1. Remove the hidden _GoBack bookmark from the paragraph that currently reads
#    "The code is like that of ref-step, ..." -- it will be re-created at the end
#    of the new content we add, matching the target document.
# 2. Replace that paragraph's text with the new, longer introduction text.
# 3. Insert a whole block of new paragraphs (explaining gpib_data.py / inst_*.py,
#    the "Bugs" heading, etc.) right after it, ending with the relocated
#    _GoBack bookmark.
# 4. Remove the old trailing "All commands are built into..." paragraph and the
#    empty paragraph that followed it.

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- locate (by index) the paragraph to rewrite ---------------------------------------
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "The code is like that of ref-step*") {
        $targetIndex = $i
        break
    }
}

# --- 1. drop the bookmark that currently sits on this paragraph -----------------------
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

# --- 2. rewrite the paragraph text ------------------------------------------------------
$d.Paragraphs.Item($targetIndex).Range.Find.Execute(
    "The code is like that of ref-step, but has no analysis at all and many things have been stripped down.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The code is similar to that of ref-step, but uses specific instrument classes. The main files that a user might want to edit are the gpib_data.py file, and all “inst_” files. Gpib_data contains the thread that communicates with all instruments and the wx table, so here the order of commands can be changed or additional commands can be added. For example, if the command to reset the lock in needs to be sent the user should add:",
    2) | Out-Null

# --- 3. insert all of the new paragraphs right after it ----------------------------------
$d.Paragraphs.Item($targetIndex).Range.InsertParagraphAfter()
$insertIndex = $targetIndex + 1

$newXml = @"
<w:p $wNs><w:r><w:t>Self.com(self.lockin.reset)</w:t></w:r></w:p>
<w:p $wNs><w:r><w:t xml:space="preserve">The self.com function executes the command in the bracket (self.lockin.reset), and also receives info on whether or not the command was successful. If it was not successful, it will call the MakeSafe on the instruments and the program stops. </w:t></w:r></w:p>
<w:p $wNs><w:r><w:t>If a new command needs to be added to the lockin, the file “inst_lockin.py” needs to be edited, a new function should be added as such:</w:t></w:r></w:p>
<w:p $wNs><w:r><w:t>def new_function(self, optional_argument):</w:t></w:r></w:p>
<w:p $wNs><w:r><w:tab/><w:t>return self.send(“The new words that need to be sent”+str(optional_argument))</w:t></w:r></w:p>
<w:p $wNs><w:r><w:t>As the name suggests, the optional argument is optinoal</w:t></w:r><w:r><w:t xml:space="preserve">. </w:t></w:r><w:r><w:t>I</w:t></w:r><w:r><w:t>f a function with an argument is called using the “self.com” method (that checks the commands sent safely), the arguments need to be given to the com function too as follows: self.com(self.some_isntrument.function, argumetns).</w:t></w:r></w:p>
<w:p $wNs><w:r><w:t>Note that sending an instrument the word ‘None’, will result in the instrument skipping the command but still printing to the event reports box that a command was skipped.</w:t></w:r></w:p>
<w:p $wNs><w:r><w:t>The instruments’ “send” function is defined in instrument.py, the instruments are all subclasses of that general class.</w:t></w:r></w:p>
<w:p $wNs/>
<w:p $wNs><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>Bugs</w:t></w:r></w:p>
<w:p $wNs><w:r><w:t>Like in ref-step, when large tables are loaded to the grid, the grid ends up creating its own scroll bar. When you mouse over that new scroll bar the computer crashes and the screen turns temporarily black. This can be avoided with smaller tables, but I don’t have a fix at the moment.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
"@

$d.Paragraphs.Item($insertIndex).Range.InsertXML($newXml)

# --- 4. remove the old trailing paragraphs ("All commands..." + the blank one) -----------
$oldIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "All commands are built into specific instrument class objects*") {
        $oldIndex = $i
        break
    }
}
$p1 = $d.Paragraphs.Item($oldIndex)
$p2 = $d.Paragraphs.Item($oldIndex + 1)
$delRange = $d.Range($p1.Range.Start, $p2.Range.End)
$delRange.Delete()
